$wb = $excel.ActiveWorkbook

# --- Rebuild the "sauceLoginsTest" sheet with new users, moving it to the end ---
# The sheet currently sits at position 2 (between sauceLoginTest and InformationPageTest).
# We Copy it to the end of the workbook (this hands the copy a fresh sheetId and lets
# later AutoFit calculations persist, unlike a brand-new Worksheets.Add() sheet), delete
# the original, then rename + refill the copy. InformationPageTest naturally shifts into
# position 2 once the original slot is vacated.
$oldSheet = $wb.Worksheets.Item("sauceLoginsTest")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldSheet.Copy($null, $lastSheet)
$oldSheet.Delete()

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "sauceLoginsTest"
$newSheet.Cells.Clear()

# --- Header + user rows ---
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("B1").Value = "Password"

$newSheet.Range("A2").Value = "standard_user"
$newSheet.Range("B2").Value = "secret_sauce"

$newSheet.Range("A3").Value = "locked_out_user"
$newSheet.Range("B3").Value = "secret_sauce"

$newSheet.Range("A4").Value = "problem_user"
$newSheet.Range("B4").Value = "secret_sauce"

$newSheet.Range("A5").Value = "performance_glitch_user"
$newSheet.Range("B5").Value = "secret_sauce"

$newSheet.Range("A6").Value = "error_user"
$newSheet.Range("B6").Value = "secret_sauce"

$newSheet.Range("A7").Value = "visual_user"
$newSheet.Range("B7").Value = "secret_sauce"

# --- Force black font colour on the used range (creates the new style) ---
$newSheet.Range("A1:B7").Font.Color = 0

# --- Column widths: best-fit to content, matching the original sheet's sizing ---
$newSheet.Columns.Item("A:B").AutoFit()

$newSheet.Activate()
